# Update "想去人数" (F column) counts on both the "展览" sheet and the
# "全部类型" sheet, which duplicate the same event rows.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 14290
$ws1.Range("F7").Value = 16325
$ws1.Range("F24").Value = 6573
$ws1.Range("F29").Value = 5708

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 14290
$ws4.Range("F7").Value = 16325
$ws4.Range("F25").Value = 6573
$ws4.Range("F32").Value = 5708
